$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       AdaBoostRegressor())]),`n                                            param_grid={'model__learning_rate': [0.1,`n                                                                                 0.5,`n                                                                                 1.0],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

# Add header for new column F (copy formatting from an existing header cell)
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update numeric values in B2:D5
$ws.Range("B2").Value = 0.4802115233506354
$ws.Range("C2").Value = 0.9904379687826885
$ws.Range("D2").Value = 0.5634498732813795

$ws.Range("B3").Value = 0.2391564293541118
$ws.Range("C3").Value = 0.995325681573941
$ws.Range("D3").Value = 0.3815729935562729

$ws.Range("B4").Value = 0.317720178047969
$ws.Range("C4").Value = 0.9938877702701309
$ws.Range("D4").Value = 0.4478343651291319

$ws.Range("B5").Value = 0.4192265508146497
$ws.Range("C5").Value = 0.9917334263902583
$ws.Range("D5").Value = 0.5047277010375186

# Add model text to F2:F5
$ws.Range("F2").Value = $modelText
$ws.Range("F3").Value = $modelText
$ws.Range("F4").Value = $modelText
$ws.Range("F5").Value = $modelText
